$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 83.59375
$ws.Range("C2").Value = 82.55208333333334
$ws.Range("D2").Value = 69.79166666666667
